$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "5-jun"
$ws.Range("D1").Value = "6-jun"
$ws.Range("E1").Value = "7-jun"
$ws.Range("F1").Value = "8-jun"

$ws.Range("C1:F1").NumberFormat = "@"

$ws.Range("F2").Select()
